$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price observation per row. A new week of
# prices (serial 44551 = 2021-12-21) is being added as three new rows at
# the top of the most recent block (rows 18-20), pushing the previously
# most-recent rows (old rows 18-21) down by three rows (to rows 21-24).

# Insert three new blank rows at row 18; this shifts the old rows 18-21
# down to 21-24 automatically, carrying their values and formatting with
# them.
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# Fill the three freshly inserted rows with the new week's entries.
# Row 18: "Especial" quality, $/caja 18 kilos boxes.
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C18").Value = "Los Lagos"
$ws.Range("D18").Value = 44551
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100103
$ws.Range("H18").Value = "Frutos de hueso (carozo)"
$ws.Range("I18").Value = 100103003
$ws.Range("J18").Value = "Damasco"
$ws.Range("K18").Value = "Castle Brite"
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 1111
$ws.Range("T18").Value = 18

# Row 19: "Primera" quality, $/caja 18 kilos boxes.
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44551
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103003
$ws.Range("J19").Value = "Damasco"
$ws.Range("K19").Value = "Castle Brite"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = "$/caja 18 kilos"
$ws.Range("R19").Value = "Región Metropolitana"
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18

# Row 20: "Segunda" quality, $/caja 18 kilos boxes.
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44551
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103003
$ws.Range("J20").Value = "Damasco"
$ws.Range("K20").Value = "Castle Brite"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 16000
$ws.Range("Q20").Value = "$/caja 18 kilos"
$ws.Range("R20").Value = "Región Metropolitana"
$ws.Range("S20").Value = 889
$ws.Range("T20").Value = 18

# Rows 21-24 now hold the former rows 18-21, already shifted down intact
# by the insert above - no further changes needed there.
